$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-01-18 Thursday", $true, $true, $false, $false, $false, $true, 1, $false, "2024-01-19 Friday", 2) | Out-Null
$d.Content.Find.Execute("218×5=", $true, $true, $false, $false, $false, $true, 1, $false, "766×6=", 2) | Out-Null
$d.Content.Find.Execute("995×8=", $true, $true, $false, $false, $false, $true, 1, $false, "977×8=", 2) | Out-Null
$d.Content.Find.Execute("354×2=", $true, $true, $false, $false, $false, $true, 1, $false, "948×8=", 2) | Out-Null
$d.Content.Find.Execute("762×3=", $true, $true, $false, $false, $false, $true, 1, $false, "330×4=", 2) | Out-Null
$d.Content.Find.Execute("417×7=", $true, $true, $false, $false, $false, $true, 1, $false, "967×5=", 2) | Out-Null
$d.Content.Find.Execute("669×6=", $true, $true, $false, $false, $false, $true, 1, $false, "710×7=", 2) | Out-Null
$d.Content.Find.Execute("294×6=", $true, $true, $false, $false, $false, $true, 1, $false, "440×3=", 2) | Out-Null
$d.Content.Find.Execute("375×9=", $true, $true, $false, $false, $false, $true, 1, $false, "130×6=", 2) | Out-Null
$d.Content.Find.Execute("691×4=", $true, $true, $false, $false, $false, $true, 1, $false, "520×6=", 2) | Out-Null
$d.Content.Find.Execute("169×5=", $true, $true, $false, $false, $false, $true, 1, $false, "360×4=", 2) | Out-Null
$d.Content.Find.Execute("963×9=", $true, $true, $false, $false, $false, $true, 1, $false, "445×4=", 2) | Out-Null
$d.Content.Find.Execute("829×4=", $true, $true, $false, $false, $false, $true, 1, $false, "543×5=", 2) | Out-Null
$d.Content.Find.Execute("718×5=", $true, $true, $false, $false, $false, $true, 1, $false, "722×5=", 2) | Out-Null
$d.Content.Find.Execute("381×3=", $true, $true, $false, $false, $false, $true, 1, $false, "618×8=", 2) | Out-Null
$d.Content.Find.Execute("605×6=", $true, $true, $false, $false, $false, $true, 1, $false, "169×3=", 2) | Out-Null
$d.Content.Find.Execute("185×3=", $true, $true, $false, $false, $false, $true, 1, $false, "821×8=", 2) | Out-Null
$d.Content.Find.Execute("672×8=", $true, $true, $false, $false, $false, $true, 1, $false, "939×2=", 2) | Out-Null
$d.Content.Find.Execute("438×9=", $true, $true, $false, $false, $false, $true, 1, $false, "871×3=", 2) | Out-Null
$d.Content.Find.Execute("660×6=", $true, $true, $false, $false, $false, $true, 1, $false, "212×2=", 2) | Out-Null
$d.Content.Find.Execute("251×6=", $true, $true, $false, $false, $false, $true, 1, $false, "163×4=", 2) | Out-Null
$d.Content.Find.Execute("186×8=", $true, $true, $false, $false, $false, $true, 1, $false, "873×4=", 2) | Out-Null
$d.Content.Find.Execute("379×5=", $true, $true, $false, $false, $false, $true, 1, $false, "908×5=", 2) | Out-Null
$d.Content.Find.Execute("354×3=", $true, $true, $false, $false, $false, $true, 1, $false, "610×5=", 2) | Out-Null
$d.Content.Find.Execute("912×5=", $true, $true, $false, $false, $false, $true, 1, $false, "805×9=", 2) | Out-Null
$d.Content.Find.Execute("106×9=", $true, $true, $false, $false, $false, $true, 1, $false, "798×5=", 2) | Out-Null
